$wb = $excel.ActiveWorkbook

$wsAreas = $wb.Worksheets.Item("areanames")
$wsDist  = $wb.Worksheets.Item("distances_fixed")

# --- Update the descriptive area-name text in column B (rows 3-12) ---
# and apply a plain black "Aptos" font to each, cell-by-cell so every
# cell converges on the same shared style entry.
$descriptions = @(
  "North America (NA) - includes Mexico except southernmost states and Yucatan",
  "Central America and Caribbean (CA) - includes southernmost Mexican states and Yucatan",
  "South America (SA) - everything south of Panama",
  "Africa (AF) - includes Arabian peninsula",
  "Madagascar and Indian Ocean islands (MA) - basically Madagascar, Mauritius, and Reunion",
  "West Eurasia (WE) - Eurasia west of Urals and Caspian Sea",
  "East Eurasia (EE) - Eurasia east of Urals and Caspian Sea, north of southern China",
  "India (IN) - Indian subcontinent, includes Pakistan, Nepal, Bangladesh",
  "Southeast Asia (SE) - basically Indomalayan biogeographic realm minus India, west of Wallace's Line",
  "Australasia (AU) - Australia, New Guinea, Wallacea, New Zealand (includes some Pacific islands close by, like Solomon Islands, Vanuatu, New Caledonia)"
)

for ($i = 0; $i -lt $descriptions.Count; $i++) {
    $row = 3 + $i
    $cell = $wsAreas.Cells.Item($row, 2)
    $cell.Value = $descriptions[$i]
    $cell.Font.Name = "Aptos"
    $cell.Font.Color = 0
}

# --- Widen column B to fit the new descriptive text ---
$wsAreas.Columns.Item(2).ColumnWidth = 133.66666666666666

# --- Update view/selection state ---
$wsDist.Activate()
$wsDist.Range("L3:U12").Select()

$wsAreas.Activate()
$wsAreas.Range("B3:B12").Select()
